# edit.ps1 -- reproduces the ITDA.docx edit (proof-reading pass: spell/grammar
# marks added around several words, plus a few genuine wording fixes).
#
# Strategy: the Word OM exposed here supports Range.InsertXML with a
# pkg:package-wrapped WordProcessingML fragment. When the Range passed to
# InsertXML spans an entire paragraph (start of paragraph through, and
# including, its paragraph mark) the call cleanly replaces that paragraph's
# contents with the <w:p> supplied in the fragment. We use that primitive to
# rebuild each touched paragraph with the extra <w:proofErr>/<w:bookmarkStart>/
# <w:bookmarkEnd> markers and run-splits the diff calls for, while leaving
# every other paragraph (and each paragraph's own pPr/rPr formatting)
# untouched.

$d = $word.ActiveDocument

$pkgOpen = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgClose = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

function Find-ParagraphByText($doc, $marker) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Contains($marker)) {
            return $p
        }
    }
    return $null
}

function Set-ParagraphXml($doc, $marker, $newParaInnerXml) {
    $p = Find-ParagraphByText $doc $marker
    if ($null -eq $p) {
        throw "paragraph not found: $marker"
    }
    $r = $p.Range
    $frag = $pkgOpen + $newParaInnerXml + $pkgClose
    $null = $r.InsertXML($frag)
}

# ---------------------------------------------------------------------
# 1. "****** Architechture *******"
# ---------------------------------------------------------------------
Set-ParagraphXml $d "Architechture" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='000C67A9'>" +
      "<w:r><w:t xml:space='preserve'>****** </w:t></w:r>" +
      "<w:proofErr w:type='spellStart'/>" +
      "<w:r><w:t>Architechture</w:t></w:r>" +
      "<w:proofErr w:type='spellEnd'/>" +
      "<w:r><w:t xml:space='preserve'> *******</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 2. "The following diagram illustrates the overall architecture of the
#     BMC TrueSight IT Data Analytics product."
# ---------------------------------------------------------------------
$rPr2 = "<w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:color w:val='000000'/><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr>"
Set-ParagraphXml $d "The following diagram illustrates" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='000C67A9'>" +
      "<w:pPr><w:pStyle w:val='NormalWeb'/><w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/><w:spacing w:before='0' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0'/><w:textAlignment w:val='top'/>$rPr2</w:pPr>" +
      "<w:r>$rPr2<w:t xml:space='preserve'>The following diagram illustrates the overall architecture of the BMC </w:t></w:r>" +
      "<w:proofErr w:type='spellStart'/>" +
      "<w:r>$rPr2<w:t>TrueSight</w:t></w:r>" +
      "<w:proofErr w:type='spellEnd'/>" +
      "<w:r>$rPr2<w:t xml:space='preserve'> IT Data Analytics product.</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 3. "This is the administrative console for the BMC TrueSight IT Data
#     Analytics product."
# ---------------------------------------------------------------------
Set-ParagraphXml $d "This is the administrative console" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='005D7CBF'>" +
      "<w:pPr><w:pStyle w:val='NormalWeb'/><w:spacing w:before='0' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0'/></w:pPr>" +
      "<w:r><w:t xml:space='preserve'>This is the administrative console for the BMC </w:t></w:r>" +
      "<w:proofErr w:type='spellStart'/>" +
      "<w:r><w:t>TrueSight</w:t></w:r>" +
      "<w:proofErr w:type='spellEnd'/>" +
      "<w:r><w:t xml:space='preserve'> IT Data Analytics product.</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 4. "The Console Server component serves as the user interface ..."
#    (note: original text has a NBSP between "BMC" and "TrueSight" --
#     preserved verbatim as &#160;)
# ---------------------------------------------------------------------
Set-ParagraphXml $d "The Console Server component serves" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='005D7CBF'>" +
      "<w:pPr><w:pStyle w:val='NormalWeb'/><w:spacing w:before='150' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0'/></w:pPr>" +
      "<w:r><w:t xml:space='preserve'>The Console Server component serves as the user interface for performing all actions in the BMC&#160;</w:t></w:r>" +
      "<w:proofErr w:type='spellStart'/>" +
      "<w:r><w:t>TrueSight</w:t></w:r>" +
      "<w:proofErr w:type='spellEnd'/>" +
      "<w:r><w:t xml:space='preserve'> IT Data Analytics product. The Console Server also acts as the gatekeeper for all actions that can be performed by using the CLI. The authentication and authorization checks for all operations are performed by the Console Server.</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 5. "The Console Server stores all the configuration details ..."
# ---------------------------------------------------------------------
Set-ParagraphXml $d "The Console Server stores all the configuration" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='005D7CBF'>" +
      "<w:pPr><w:pStyle w:val='NormalWeb'/><w:spacing w:before='150' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0'/></w:pPr>" +
      "<w:r><w:t xml:space='preserve'>The Console Server stores all the configuration details in the Configuration Database. </w:t></w:r>" +
      "<w:proofErr w:type='gramStart'/>" +
      "<w:r><w:t>This components</w:t></w:r>" +
      "<w:proofErr w:type='gramEnd'/>" +
      "<w:r><w:t xml:space='preserve'> interacts with all the other product components.</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 6. "The Indexer component acts as an internal datastore ..."
# ---------------------------------------------------------------------
$rPr6 = "<w:rPr><w:color w:val='000000'/></w:rPr>"
Set-ParagraphXml $d "The Indexer component acts as an internal datastore" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='005D7CBF'>" +
      "<w:pPr><w:pStyle w:val='NormalWeb'/><w:spacing w:before='0' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0'/></w:pPr>" +
      "<w:r>$rPr6<w:t xml:space='preserve'>The Indexer component acts as an internal datastore used by BMC </w:t></w:r>" +
      "<w:proofErr w:type='spellStart'/>" +
      "<w:r>$rPr6<w:t>TrueSight</w:t></w:r>" +
      "<w:proofErr w:type='spellEnd'/>" +
      "<w:r>$rPr6<w:t xml:space='preserve'> IT Data Analytics, for storing all the data that is collected by using Collection Stations and Collection Agents.</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 7. "All configurations of data collectors, views, and saved searches..."
# ---------------------------------------------------------------------
Set-ParagraphXml $d "All configurations of data collectors" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='005D7CBF'>" +
      "<w:pPr><w:pStyle w:val='NormalWeb'/><w:spacing w:before='0' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0'/></w:pPr>" +
      "<w:r><w:t xml:space='preserve'>All configurations of data collectors, views, and saved searches that are used in BMC </w:t></w:r>" +
      "<w:proofErr w:type='spellStart'/>" +
      "<w:r><w:t>TrueSight</w:t></w:r>" +
      "<w:proofErr w:type='spellEnd'/>" +
      "<w:r><w:t xml:space='preserve'> IT Data Analytics are stored in the Configuration Database.</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 8. "Speed up MTTR by using advanced analytics with your log data"
#    -> "Speed up MTTR (mean time to repair) by using advanced analytics..."
# ---------------------------------------------------------------------
Set-ParagraphXml $d "Speed up MTTR" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='000C67A9'>" +
      "<w:r><w:t>Speed up MTTR</w:t></w:r>" +
      "<w:r><w:t xml:space='preserve'> (mean time to repair)</w:t></w:r>" +
      "<w:r><w:t xml:space='preserve'> by using advanced analytics with your log data</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 9. "-Focus on whats not normal by comparing to a known normal(good)"
# ---------------------------------------------------------------------
Set-ParagraphXml $d "-Focus on whats not normal" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='000C67A9'>" +
      "<w:r><w:t xml:space='preserve'>-Focus on </w:t></w:r>" +
      "<w:proofErr w:type='spellStart'/>" +
      "<w:r><w:t>whats</w:t></w:r>" +
      "<w:proofErr w:type='spellEnd'/>" +
      "<w:r><w:t xml:space='preserve'> not normal by comparing to a known normal(good)</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 10a. "-Takes coalescs to the next level" -> "-Takes coalesces to the next
#       level" with the _GoBack bookmark relocated here.
# ---------------------------------------------------------------------
Set-ParagraphXml $d "-Takes coalescs to the next level" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='000C67A9'>" +
      "<w:r><w:t xml:space='preserve'>-Takes </w:t></w:r>" +
      "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
      "<w:r><w:t>coalesces</w:t></w:r>" +
      "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
      "<w:bookmarkEnd w:id='0'/>" +
      "<w:r><w:t>to the next level</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 10b. "-Is there anything unsual in the log data? More of a particular
#       type of mesage? Less?"
# ---------------------------------------------------------------------
Set-ParagraphXml $d "-Is there anything unsual" (
    "<w:p w:rsidR='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='000C67A9'>" +
      "<w:r><w:t xml:space='preserve'>-Is there anything </w:t></w:r>" +
      "<w:proofErr w:type='spellStart'/>" +
      "<w:r><w:t>unsual</w:t></w:r>" +
      "<w:proofErr w:type='spellEnd'/>" +
      "<w:r><w:t xml:space='preserve'> in the log data? More of a </w:t></w:r>" +
      "<w:proofErr w:type='gramStart'/>" +
      "<w:r><w:t xml:space='preserve'>particular type of </w:t></w:r>" +
      "<w:proofErr w:type='spellStart'/>" +
      "<w:r><w:t>mesage</w:t></w:r>" +
      "<w:proofErr w:type='spellEnd'/>" +
      "<w:proofErr w:type='gramEnd'/>" +
      "<w:r><w:t>? Less?</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------
# 11. Remove the old _GoBack bookmark that used to sit after "Upload File"
#     (it moved to step 10a above).
# ---------------------------------------------------------------------
Set-ParagraphXml $d "Upload File" (
    "<w:p w:rsidR='000C67A9' w:rsidRPr='000C67A9' w:rsidRDefault='000C67A9' w:rsidP='000C67A9'>" +
      "<w:r><w:t>Upload File</w:t></w:r>" +
      "<w:r><w:rPr><w:rFonts w:ascii='Consolas' w:eastAsia='Times New Roman' w:hAnsi='Consolas' w:cs='Times New Roman'/><w:color w:val='000000'/><w:sz w:val='21'/><w:szCs w:val='21'/><w:lang w:eastAsia='en-IN'/></w:rPr><w:br/></w:r>" +
    "</w:p>"
)

Write-Host "done"
